$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data
# D-column values that are valid numbers need NumberFormat "@" forced
# so Excel keeps them as text (matching the original inline-string cells)
# instead of auto-converting them into numeric values.

# Row 2
$ws.Range("D2").Value = "68.273.35"
$ws.Range("E2").Value = "  +1.42%  "

# Row 3
$ws.Range("D3").Value = "3.563.63"
$ws.Range("E3").Value = "  +2.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.12"
$ws.Range("E5").Value = "  +2.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.27"
$ws.Range("E6").Value = "  +4.43%  "

# Row 7
$ws.Range("D7").Value = "3.561.19"
$ws.Range("E7").Value = "  +1.98%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  +2.04%  "

# Row 10
$ws.Range("E10").Value = "  +5.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.33"
$ws.Range("E11").Value = "  +5.17%  "

# Row 12
$ws.Range("E12").Value = "  +3.65%  "

# Row 13
$ws.Range("E13").Value = "  +2.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.19"
$ws.Range("E14").Value = "  +5.52%  "

# Row 15
$ws.Range("D15").Value = "4.165.05"
$ws.Range("E15").Value = "  +2.05%  "

# Row 16
$ws.Range("D16").Value = "3.566.29"
$ws.Range("E16").Value = "  +2.04%  "

# Row 17
$ws.Range("D17").Value = "68.199.93"
$ws.Range("E17").Value = "  +1.46%  "

# Row 18
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
$ws.Range("E19").Value = "  +6.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.03"
$ws.Range("E20").Value = "  +6.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.00"
$ws.Range("E21").Value = "  +10.97%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.69"
$ws.Range("E22").Value = "  +1.85%  "

# Row 23
$ws.Range("E23").Value = "  +3.53%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.57"
$ws.Range("E24").Value = "  +1.95%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000130"
$ws.Range("E25").Value = "  +0.98%  "

# Row 26
$ws.Range("D26").Value = "3.702.55"
$ws.Range("E26").Value = "  +2.00%  "

# Row 27
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.56"
$ws.Range("E28").Value = "  +4.34%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.08"
$ws.Range("E29").Value = "  +9.42%  "

# Row 30
$ws.Range("E30").Value = "  +3.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.70"
$ws.Range("E31").Value = "  +8.86%  "

# Row 32
$ws.Range("E32").Value = "  +6.18%  "

# Row 33
$ws.Range("E33").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.42"
$ws.Range("E34").Value = "  +4.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.13"
$ws.Range("E35").Value = "  +1.79%  "

# Row 36
$ws.Range("E36").Value = "  +3.52%  "

# Row 37
$ws.Range("D37").Value = "3.556.25"
$ws.Range("E37").Value = "  +2.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.28"
$ws.Range("E38").Value = "  +3.29%  "

# Row 39
$ws.Range("E39").Value = "  +7.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "178.84"
$ws.Range("E41").Value = "  +3.10%  "

# Row 42
$ws.Range("E42").Value = "  +4.90%  "

# Row 43
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.59"
$ws.Range("E44").Value = "  +3.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.02"
$ws.Range("E45").Value = "  +13.90%  "

# Row 46
$ws.Range("E46").Value = "  +1.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.54"
$ws.Range("E47").Value = "  +2.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.34"
$ws.Range("E48").Value = "  +6.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.67"
$ws.Range("E49").Value = "  +3.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.80"
$ws.Range("E50").Value = "  +3.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.261"
$ws.Range("E51").Value = "  +7.14%  "
